# Weekly update: a new "Mora" price observation was recorded for
# Mercado Mayorista Lo Valledor de Santiago (Provincia de Curicó).
# It belongs chronologically right after the existing header/ramp-up
# rows, so it is inserted as the new row 26, pushing the previous
# rows 26-49 down to 27-50 (the sheet's used range grows from
# A1:T49 to A1:T50).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 26 - this shifts existing rows
# 26:49 down to 27:50 and keeps their values/formatting intact.
$ws.Rows("26:26").Insert()

# Populate the newly-inserted row 26 with the new observation.
$ws.Range("A26").Value = 6
$ws.Range("B26").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C26").Value = "Metropolitana"
$ws.Range("D26").Value2 = 44546
$ws.Range("E26").Value = 13
$ws.Range("F26").Value = "Fruta"
$ws.Range("G26").Value = 100101
$ws.Range("H26").Value = "Berries"
$ws.Range("I26").Value = 100101008
$ws.Range("J26").Value = "Mora"
$ws.Range("K26").Value = "Sin especificar"
$ws.Range("L26").Value = "Primera"
$ws.Range("M26").Value = 250
$ws.Range("N26").Value = 5000
$ws.Range("O26").Value = 5000
$ws.Range("P26").Value = 5000
$ws.Range("Q26").Value = "$/bandeja 2 kilos"
$ws.Range("R26").Value = "Provincia de Curicó"
$ws.Range("S26").Value = 2500
$ws.Range("T26").Value = 2
